$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (columns C,D,E,F,G get narrower; J,K are new) ---
$ws.Columns.Item(3).ColumnWidth = 18.666666666666668
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Columns.Item(5).ColumnWidth = 13.5
$ws.Columns.Item(6).ColumnWidth = 13.666666666666666
$ws.Columns.Item(7).ColumnWidth = 13.666666666666666
$ws.Columns.Item(10).ColumnWidth = 13.0
$ws.Columns.Item(11).ColumnWidth = 15.166666666666666

# --- Header row: add new "expected"/"keterangan" columns ---
$ws.Range("J1").Value = "expected"
$ws.Range("K1").Value = "keterangan"

# --- Row 2 (existing valid "ducati" record): mark as passed ---
$ws.Range("J2").Value = "passed"

# --- Row 3 used to hold the ferrari/spider test record; that data moves to
#     row 4 with deliberately-wrong values, so clear the old row 3 fields
#     and use row 3 only for the "failed"/"fieldempty" note of the new case ---
$ws.Range("A3:I3").ClearContents()
$ws.Range("J3").Value = "failed"
$ws.Range("K3").Value = "fieldempty"

# --- Row 4: new test case with intentionally malformed values ---
$ws.Range("A4").Value = "ferrari "
$ws.Range("B4").Value = "monster"
$ws.Range("C4").Value = "795(CKD)"
$ws.Range("D4").Value = 2015
$ws.Range("E4").Value = 2017
$ws.Range("H4").Value = "automatic"
$ws.Range("I4").Value = "terbuka"
$ws.Range("J4").Value = "failed"
$ws.Range("K4").Value = "brandsalah"

# --- Final selection as left by the author ---
$ws.Range("B5").Select()
